$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# New order rows appended below the existing data (rows 2-61).
# Values in columns A and F are numeric-looking text, so they are written
# with a leading apostrophe to keep them stored as text (matching the rest
# of the sheet, which stores every value - including numbers - as text).
$ws.Cells.Item(62, 3).Value  = "418_松虫草白_scabiosa white_undefined_1bunch"
$ws.Cells.Item(62, 6).Value  = "'5"

$ws.Cells.Item(63, 1).Value  = "'10"
$ws.Cells.Item(63, 3).Value  = "649_洋牡丹樱花粉_undefined_undefined_1bunch"
$ws.Cells.Item(63, 6).Value  = "'10"

$ws.Cells.Item(64, 3).Value  = "718_银莲白_undefined_undefined_1bunch"
$ws.Cells.Item(64, 6).Value  = "'5"

$ws.Cells.Item(65, 3).Value  = "691_银莲紫_undefined_undefined_1bunch"
$ws.Cells.Item(65, 6).Value  = "'5"

$ws.Cells.Item(66, 3).Value  = "681_锦鲤橙_undefined_undefined_1bunch"
$ws.Cells.Item(66, 6).Value  = "'5"

$ws.Cells.Item(67, 3).Value  = "587_洋牡丹橙_undefined_undefined_1bunch"
$ws.Cells.Item(67, 6).Value  = "'10"

$ws.Cells.Item(68, 3).Value  = "585_洋牡丹红_undefined_undefined_1bunch"
$ws.Cells.Item(68, 6).Value  = "'15"

$ws.Cells.Item(69, 3).Value  = "100_绣球单瓣白_Hydrangea White S_Hydrangea L._1stem"
$ws.Cells.Item(69, 6).Value  = "'11"

$ws.Cells.Item(70, 1).Value  = "'1"

# Update the Summary sheet's encoded Number-column digest (G2) to include
# the eight newly-added "Number" values (plus the trailing blank row 70).
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Cells.Item(2, 7).Value = "'02010201055555301051510301030151515151210101555101891510205712881191041351155302010205310102020155105105551015110"
